# Updates Leve profit calculation columns (currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed Market Board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3904.3572
$ws.Range("J43").Value = 2108.6
$ws.Range("L43").Value = 2108.6
$ws.Range("N43").Value = -2246.6

$ws.Range("H107").Value = 889.4
$ws.Range("I107").Value = 889.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 889.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1030.6
$ws.Range("N107").ClearContents()

$ws.Range("H138").Value = 297995.62
$ws.Range("I138").Value = 3793.0303
$ws.Range("J138").Value = 484701.12
$ws.Range("K138").Value = 11379.0909
$ws.Range("L138").Value = 1454103.36
$ws.Range("M138").Value = -6239.090899999999
$ws.Range("N138").Value = -1464383.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1125
$ws.Range("I4").Value = 750
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 750
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -634
$ws.Range("N4").Value = -1732

$ws.Range("H32").Value = 5415.846
$ws.Range("I32").Value = 3991.9138
$ws.Range("K32").Value = 3991.9138
$ws.Range("M32").Value = -3704.9138

$ws.Range("H61").Value = 9476.25
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 9476.25
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 9476.25
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -9900.25

$ws.Range("H132").Value = 2981.318
$ws.Range("I132").Value = 2513.5
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 7540.5
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -5010.5
$ws.Range("N132").Value = -16460

$ws.Range("H136").Value = 9476.25
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 9476.25
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 28428.75
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -33528.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I105").Value = 386513.84
$ws.Range("K105").Value = 386513.84
$ws.Range("M105").Value = -384766.84

$ws.Range("H107").Value = 5129709.5
$ws.Range("I107").Value = 5918470
$ws.Range("K107").Value = 5918470
$ws.Range("M107").Value = -5916550

$ws.Range("H134").Value = 2443.7673
$ws.Range("I134").Value = 2089.9412
$ws.Range("K134").Value = 6269.823600000001
$ws.Range("M134").Value = -3734.823600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 76927720
$ws.Range("I7").Value = 6920
$ws.Range("K7").Value = 6920
$ws.Range("M7").Value = -6807

$ws.Range("H31").Value = 4212.8696
$ws.Range("I31").Value = 4151.5
$ws.Range("K31").Value = 4151.5
$ws.Range("M31").Value = -3856.5

$ws.Range("H34").Value = 4212.8696
$ws.Range("I34").Value = 4151.5
$ws.Range("K34").Value = 4151.5
$ws.Range("M34").Value = -3949.5

$ws.Range("H58").Value = 4761.1333
$ws.Range("I58").Value = 5999.5
$ws.Range("J58").Value = 4570.615
$ws.Range("K58").Value = 5999.5
$ws.Range("L58").Value = 4570.615
$ws.Range("M58").Value = -5796.5
$ws.Range("N58").Value = -4976.615

$ws.Range("H86").Value = 4110.8184
$ws.Range("I86").Value = 3965.375
$ws.Range("J86").Value = 4498.6665
$ws.Range("K86").Value = 3965.375
$ws.Range("L86").Value = 4498.6665
$ws.Range("M86").Value = -2842.375
$ws.Range("N86").Value = -6744.6665

$ws.Range("H89").Value = 4110.8184
$ws.Range("I89").Value = 3965.375
$ws.Range("J89").Value = 4498.6665
$ws.Range("K89").Value = 19826.875
$ws.Range("L89").Value = 22493.3325
$ws.Range("M89").Value = -14210.875
$ws.Range("N89").Value = -33725.3325

$ws.Range("H115").Value = 49907.5
$ws.Range("J115").Value = 49907.5
$ws.Range("L115").Value = 49907.5
$ws.Range("N115").Value = -52257.5

$ws.Range("H132").Value = 31254758
$ws.Range("I132").Value = 50004300
$ws.Range("K132").Value = 150012900
$ws.Range("M132").Value = -150010370

$ws.Range("H134").Value = 3097.6
$ws.Range("I134").Value = 2809.6667
$ws.Range("K134").Value = 8429.000100000001
$ws.Range("M134").Value = -5894.000100000001

$ws.Range("H136").Value = 4761.1333
$ws.Range("I136").Value = 5999.5
$ws.Range("J136").Value = 4570.615
$ws.Range("K136").Value = 17998.5
$ws.Range("L136").Value = 13711.845
$ws.Range("M136").Value = -15448.5
$ws.Range("N136").Value = -18811.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 563.5789
$ws.Range("J117").Value = 977.8570999999999
$ws.Range("L117").Value = 2933.5713
$ws.Range("N117").Value = -9817.5713

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

$ws.Range("H129").Value = 2098.625
$ws.Range("I129").Value = 1805
$ws.Range("J129").Value = 2685.875
$ws.Range("K129").Value = 5415
$ws.Range("L129").Value = 8057.625
$ws.Range("M129").Value = -415
$ws.Range("N129").Value = -18057.625

$ws.Range("H131").Value = 7029.5
$ws.Range("I131").Value = 17297.867
$ws.Range("J131").Value = 2060.9355
$ws.Range("K131").Value = 51893.601
$ws.Range("L131").Value = 6182.806500000001
$ws.Range("M131").Value = -46853.601
$ws.Range("N131").Value = -16262.8065

$ws.Range("H137").Value = 3680.7083
$ws.Range("I137").Value = 3860.6365
$ws.Range("J137").Value = 3528.4614
$ws.Range("K137").Value = 11581.9095
$ws.Range("L137").Value = 10585.3842
$ws.Range("M137").Value = -6481.9095
$ws.Range("N137").Value = -20785.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3228.1936
$ws.Range("I97").Value = 773
$ws.Range("J97").Value = 11646
$ws.Range("K97").Value = 773
$ws.Range("L97").Value = 11646
$ws.Range("M97").Value = -277
$ws.Range("N97").Value = -12638

$ws.Range("H102").Value = 1029.0416
$ws.Range("I102").Value = 829.825
$ws.Range("K102").Value = 829.825
$ws.Range("M102").Value = 792.175

$ws.Range("H126").Value = 6643
$ws.Range("J126").Value = 8257.714
$ws.Range("L126").Value = 24773.142
$ws.Range("N126").Value = -29713.142

$ws.Range("H132").Value = 2885.5557
$ws.Range("I132").Value = 2281.4285
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6844.2855
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4314.2855
$ws.Range("N132").Value = -20060

$ws.Range("H136").Value = 10921.823
$ws.Range("J136").Value = 10921.823
$ws.Range("L136").Value = 32765.469
$ws.Range("N136").Value = -37865.469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5233.2856
$ws.Range("I7").Value = 4355.2144
$ws.Range("J7").Value = 6989.4287
$ws.Range("K7").Value = 4355.2144
$ws.Range("L7").Value = 6989.4287
$ws.Range("M7").Value = -4243.2144
$ws.Range("N7").Value = -7213.4287

$ws.Range("H68").Value = 1073.75
$ws.Range("J68").Value = 995
$ws.Range("L68").Value = 995
$ws.Range("N68").Value = -2493

$ws.Range("H71").Value = 1073.75
$ws.Range("J71").Value = 995
$ws.Range("L71").Value = 4975
$ws.Range("N71").Value = -12463

$ws.Range("H126").Value = 5233.2856
$ws.Range("I126").Value = 4355.2144
$ws.Range("J126").Value = 6989.4287
$ws.Range("K126").Value = 13065.6432
$ws.Range("L126").Value = 20968.2861
$ws.Range("M126").Value = -10595.6432
$ws.Range("N126").Value = -25908.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 995
$ws.Range("J3").Value = 995
$ws.Range("L3").Value = 995
$ws.Range("N3").Value = -1223

$ws.Range("H81").Value = 4061.762
$ws.Range("I81").Value = 3581.9
$ws.Range("J81").Value = 4498
$ws.Range("K81").Value = 7163.8
$ws.Range("L81").Value = 8996
$ws.Range("M81").Value = -6102.8
$ws.Range("N81").Value = -11118

$ws.Range("H84").Value = 4061.762
$ws.Range("I84").Value = 3581.9
$ws.Range("J84").Value = 4498
$ws.Range("K84").Value = 35819
$ws.Range("L84").Value = 44980
$ws.Range("M84").Value = -30515
$ws.Range("N84").Value = -55588

$ws.Range("H96").Value = 7839.4287
$ws.Range("I96").Value = 7820
$ws.Range("K96").Value = 7820
$ws.Range("M96").Value = -6447

$ws.Range("H132").Value = 7578129
$ws.Range("I132").Value = 8549315
$ws.Range("K132").Value = 25647945
$ws.Range("M132").Value = -25645415

$ws.Range("H136").Value = 27780506
$ws.Range("I136").Value = 38462776
$ws.Range("J136").Value = 6599.9
$ws.Range("K136").Value = 115388328
$ws.Range("L136").Value = 19799.7
$ws.Range("M136").Value = -115385778
$ws.Range("N136").Value = -24899.7
